$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Experiment 1 block: header label change ---
$ws.Range("F3").Value = "Energy at t=0.1"

# --- Experiment 1 block: new raw data values (row 4-8) ---
$ws.Range("B4").Value = 0.055277508705783802
$ws.Range("F4").Value = 0.153066905474987

$ws.Range("B5").Value = 0.013387753817850701
$ws.Range("F5").Value = 0.173677776602322

$ws.Range("B6").Value = 0.0033089979620469102
$ws.Range("F6").Value = 0.17907953104893801

$ws.Range("B7").Value = 0.00080630448057797399
$ws.Range("F7").Value = 0.18044882360572001

$ws.Range("F8").Value = 0.18088862945063999

# --- Experiment 2 block: new raw data values (rows 12-17) ---
$ws.Range("C12").Value = 0.0000042759255746980998
$ws.Range("E12").Value = -0.025427564566334101

$ws.Range("C13").Value = 0.00000106881441369226
$ws.Range("E13").Value = -0.025427212062880501

$ws.Range("C14").Value = 0.00000026718239815483
$ws.Range("E14").Value = -0.0254271239471488

$ws.Range("C15").Value = 0.000000066795472854015002
$ws.Range("E15").Value = -0.025427101920295101

$ws.Range("C16").Value = 0.000000016690242116144301
$ws.Range("E16").Value = -0.025427096412487801

$ws.Range("E17").Value = -0.025427094577657999

# --- Experiment 3 block: new raw data values (rows 21-26) ---
$ws.Range("C21").Value = 0.079777383804321297
$ws.Range("D21").Value = 0.079777431488037107

$ws.Range("C22").Value = 0.0223781786571066
$ws.Range("D22").Value = 0.0223781976775171

$ws.Range("C23").Value = 0.0077581453496418399
$ws.Range("D23").Value = 0.0077581506850154004

$ws.Range("C24").Value = 0.0019887971747154501
$ws.Range("D24").Value = 0.0019887990244014399

$ws.Range("C25").Value = 0.00040725399661228402
$ws.Range("D25").Value = 0.00040725423369538202

$ws.Range("C26").Value = 0.000091693572748853503
$ws.Range("D26").Value = 0.000091693621297313094

# --- New empty, number-formatted cells (K/L/N/O columns, rows 16-21) ---
$ws.Range("K16:L16").NumberFormat = "0.00E+00"
$ws.Range("N16:O16").NumberFormat = "0.00E+00"
$ws.Range("K17:L17").NumberFormat = "0.00E+00"
$ws.Range("N17:O17").NumberFormat = "0.00E+00"
$ws.Range("K18:L18").NumberFormat = "0.00E+00"
$ws.Range("N18:O18").NumberFormat = "0.00E+00"
$ws.Range("K19:L19").NumberFormat = "0.00E+00"
$ws.Range("N19:O19").NumberFormat = "0.00E+00"
$ws.Range("K20:L20").NumberFormat = "0.00E+00"
$ws.Range("N20:O20").NumberFormat = "0.00E+00"
$ws.Range("K21").NumberFormat = "0.00E+00"
$ws.Range("N21").NumberFormat = "0.00E+00"

# --- Sheet view: update selected cell ---
$null = $ws.Range("F33").Select()

# --- Page setup: force explicit portrait orientation ---
$ws.PageSetup.Orientation = 1
